# Apply changes described by the commit diff:
# 1. Replace "KAS" with "CAN" in the fertilizer type cells.
# 2. Clear stray leftover values (kept only formatting) in several
#    "Fertilization" sub-table rows (the 4th and subsequent fertilization
#    events no longer have date/amount/fertilizer-type data recorded).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# 1. KAS -> CAN
$kasCells = @("E47","E48","E58","E59","E60","E69","E70","E71")
foreach ($addr in $kasCells) {
    $ws.Range($addr).Value = "CAN"
}

# 2. Clear contents (keep formatting) of now-empty cells across the three
#    fertilization blocks (rows 50-53, 61-64, 72-75).
$clearCells = @(
    "B50","C50","D50","E50",
    "C51","E51",
    "C52","D52","E52",
    "C53","D53","E53",
    "B61","C61","D61","E61",
    "C62","E62",
    "C63","D63","E63",
    "C64","D64","E64",
    "B72","C72","D72","E72",
    "C73","E73",
    "C74","D74","E74",
    "C75","D75","E75"
)
foreach ($addr in $clearCells) {
    $ws.Range($addr).ClearContents()
}

# Reflect the final cursor/selection position recorded in the saved file.
$ws.Range("E71").Select() | Out-Null
